$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $searchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$searchText*") {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) Intro paragraph: "Please create 4 rows of four stylized..." ->
#    "For the highlighted tiles, please create a single row of
#     stylized hyperlink tiles, each with a vibrant and modern
#     design. Here's a detailed breakdown of the components:"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Please create 4 rows of four stylized hyperlink tiles, each with a vibrant and modern design. Here" + [char]0x2019 + "s a detailed breakdown of its components:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "For the highlighted tiles, please create a single row of stylized hyperlink tiles, each with a vibrant and modern design. Here" + [char]0x2019 + "s a detailed breakdown of the components:",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) "Row Structure" paragraph: horizontally -> vertically
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "The tiles are arranged horizontally in a single row",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The tiles are arranged vertically in a single row",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) "Colors" paragraph: add ", but professional"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "The color palette is vibrant and eye-catching, making each tile stand out distinctly.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The color palette is vibrant and eye-catching, but professional, making each tile stand out distinctly.",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) Insert two new paragraphs right after "These buttons add an
#    interactive element..." and before "Aesthetic and Design:":
#      a) a numbered (ilvl=1, numId=5) bullet about the scroll arrows
#      b) an empty paragraph indented at 1440 twips (1 inch)
# ------------------------------------------------------------------
$anchor = Find-ParagraphByText $d "These buttons add an interactive element"

$anchor.Range.InsertParagraphAfter() | Out-Null
$scrollPara = $anchor.Next()
$scrollPara.Range.Text = "To scroll left or right on the screen, create 2 white circles with simple thin blue arrows, one pointing to the left and the other one pointing to the right."

$scrollPara.Range.InsertParagraphAfter() | Out-Null
$blankPara = $scrollPara.Next()
$blankPara.Range.ListFormat.RemoveNumbers() | Out-Null
$blankPara.Format.LeftIndent = 72
